# Refactor timetable generation to apply consistent cell alignment and
# borders for improved visual presentation: center text (horizontally and
# vertically) and enable text wrapping across the whole table; also fix the
# column-G (Saturday) schedule so the call-time / concert / refreshment
# entries line up with the correct time slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Consistent alignment + wrap across the whole table -----------------
# Every populated / bordered cell in the grid (A1:G30) ends up centered both
# horizontally and vertically, with text wrapping enabled.
$all = $ws.Range("A1:G30")
$all.HorizontalAlignment = -4108   # xlCenter
$all.VerticalAlignment = -4108     # xlCenter
$all.WrapText = $true

# --- 2. Fix up the Saturday (column G) afternoon schedule -------------------
# Previously:
#   G15:G18 "Lunch / Dress Up, Warm Up"   merged
#   G19      "Concert call time"          (own cell)
#   G20:G25 "Lina Summer Camp ... Concert" merged
#   G26:G28 "After concert refreshment ..." merged
# Correct layout (shifted down one row starting at G19):
#   G15:G19 "Lunch / Dress Up, Warm Up"   merged (extended to include G19)
#   G20      "Concert call time"          (own cell)
#   G21:G26 "Lina Summer Camp ... Concert" merged
#   G27:G28 "After concert refreshment ..." merged

# Unmerge the ranges that need to change shape first.
$ws.Range("G15:G18").UnMerge()
$ws.Range("G20:G25").UnMerge()
$ws.Range("G26:G28").UnMerge()

# Move the cell contents down into their corrected slots.
$ws.Range("G19").ClearContents()
$ws.Range("G20").Value = "Concert call time"
$ws.Range("G21").Value = "Lina Summer Camp of Music Students & Friends Concert"
$ws.Range("G26").ClearContents()
$ws.Range("G27").Value = "After concert refreshment " + [char]10 + "(Maritime Museum)"

# Re-merge into the corrected ranges.
$ws.Range("G15:G19").Merge()
$ws.Range("G21:G26").Merge()
$ws.Range("G27:G28").Merge()
